# Adding BossPunching_Yellow.png 3 frames, 147x67 sheet size
# Row 9 on sheet "BOSS YELLOW" corresponds to pose #5 "Punching":
#   D9 (FILENAME)   : (blank) -> "BossPunching_Yellow.png"
#   F9 (FRAMES)     : 4 -> 3
#   G9 (SHEET SIZE) : "196x67" -> "147x67"
#   H9 (STATUS)     : (blank) -> "DONE"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "BossPunching_Yellow.png"
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = "147x67"
$ws.Range("H9").Value = "DONE"

# Update active cell selection to D10, matching the author's final cursor position.
$ws.Range("D10").Select()
